# Insert two new rows at row 17 (pushing existing rows 17.. down to 19..)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()

# Row 17: new data row (Banquete)
$ws.Cells.Item(17, 1).Value = 12
$ws.Cells.Item(17, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(17, 3).Value = "Metropolitana"
$ws.Cells.Item(17, 4).Value = 44525
$ws.Cells.Item(17, 5).Value = 13
$ws.Cells.Item(17, 6).Value = 300000000
$ws.Cells.Item(17, 7).Value = "Espárragos"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Banquete"
$ws.Cells.Item(17, 10).Value = 500
$ws.Cells.Item(17, 11).Value = 1200
$ws.Cells.Item(17, 12).Value = 1200
$ws.Cells.Item(17, 13).Value = 1200
$ws.Cells.Item(17, 14).Value = "$/kilo"
$ws.Cells.Item(17, 15).Value = "Provincia de Linares"
$ws.Cells.Item(17, 16).Value = 1200
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = "Hortaliza"

# Row 18: new data row (Primera)
$ws.Cells.Item(18, 1).Value = 12
$ws.Cells.Item(18, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(18, 3).Value = "Metropolitana"
$ws.Cells.Item(18, 4).Value = 44525
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(18, 6).Value = 300000000
$ws.Cells.Item(18, 7).Value = "Espárragos"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 440
$ws.Cells.Item(18, 11).Value = 1000
$ws.Cells.Item(18, 12).Value = 1000
$ws.Cells.Item(18, 13).Value = 1000
$ws.Cells.Item(18, 14).Value = "$/kilo"
$ws.Cells.Item(18, 15).Value = "Provincia de Linares"
$ws.Cells.Item(18, 16).Value = 1000
$ws.Cells.Item(18, 17).Value = 1
$ws.Cells.Item(18, 18).Value = "Hortaliza"
